# #5: property boat&car done
# Extend the "汽車" (car) sheet with the same metadata columns used on the
# other property sheets: property_category, category, date, legislator_name,
# legislator_id, source_file, index — and add a new "capacity" header column
# right after "name".

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("汽車")

# --- Header row (row 1) ---------------------------------------------------
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "capacity"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "register_date"
$ws.Cells.Item(1, 6).Value = "register_reason"
$ws.Cells.Item(1, 7).Value = "acquire_value"
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# --- Data row (row 2) ------------------------------------------------------
$ws.Cells.Item(2, 2).Value = "nissanlivina"
$ws.Cells.Item(2, 3).Value = 1598
$ws.Cells.Item(2, 4).Value = "林君倩"
$ws.Cells.Item(2, 5).Value = "99年05月27日"
$ws.Cells.Item(2, 6).Value = "買賣"
$ws.Cells.Item(2, 7).Value = 568000
$ws.Cells.Item(2, 8).Value = "land"
$ws.Cells.Item(2, 9).Value = "normal"

# "2012-04-20" reads like a real date to Excel's auto-detection, so force
# text formatting while assigning it (restored to the row's plain format
# below, along with the rest of the new columns).
$ws.Cells.Item(2, 10).NumberFormat = "@"
$ws.Cells.Item(2, 10).Value = "2012-04-20"

$ws.Cells.Item(2, 11).Value = "楊曜"
$ws.Cells.Item(2, 12).Value = 1759
$ws.Cells.Item(2, 13).Value = "tmpcdc61"
$ws.Cells.Item(2, 14).Value = 30

# Extend the header (row 1, bold+bordered) and data (row 2, plain) look
# across the new columns H:N by copying the existing formatted cells, so
# the new cells match the rest of the sheet (this also normalises the
# temporary text format used for the date cell above back to the plain
# data-row look).
$ws.Range("B1:G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B2:G2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0
